$wb = $excel.ActiveWorkbook

# The workbook has two sheets carrying the same "展览" (exhibition) table:
# "展览" and "全部类型". Update the "想去人数" (F column) figures on both.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1907
    $ws.Range("F4").Value = 1169
    $ws.Range("F5").Value = 1244
    $ws.Range("F7").Value = 6019
}
